# Calculations Reference.xlsx edit script
# - Adds a new blank worksheet "Juml" between "PWM" and "Cal_4-30"
# - Adds a pump speed control table (Creep, Slow, Fast, Stop) to the "PWM" sheet

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Juml" worksheet, positioned after PWM / before Cal_4-30 ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Juml"
$juml = $wb.Worksheets.Item("Juml")
$calSheet = $wb.Worksheets.Item("Cal_4-30")
$juml.Move($calSheet)

# --- 2. Populate the new pump speed control table on the "PWM" sheet ---
$pwm = $wb.Worksheets.Item("PWM")

# Header row
$pwm.Range("C26").Value = "R1"
$pwm.Range("D26").Value = "R2"
$pwm.Range("E26").Value = "Ton = R1 + R2"
$pwm.Range("F26").Value = "T = R1 +2R2"
$pwm.Range("G26").Value = "DC"
$pwm.Range("G26").Style = "Normal"

# Creep row
$pwm.Range("B27").Value = "Creep"
$pwm.Range("C27").Value = 10000
$pwm.Range("D27").Value = 2200
$pwm.Range("E27").Formula = "=C27+D27"
$pwm.Range("F27").Formula = "=C27+2*D27"
$pwm.Range("G27").Formula = "=E27/F27"
$pwm.Range("G27").Style = "Normal"

# Slow row
$pwm.Range("B28").Value = "Slow"
$pwm.Range("C28").Value = 750
$pwm.Range("D28").Value = 10000
$pwm.Range("E28").Formula = "=C28+D28"
$pwm.Range("F28").Formula = "=C28+2*D28"
$pwm.Range("G28").Formula = "=E28/F28"
$pwm.Range("G28").Style = "Normal"

# Fast row
$pwm.Range("B29").Value = "Fast"
$pwm.Range("C29").Value = 0
$pwm.Range("D29").Value = 4777
$pwm.Range("E29").Formula = "=C29+D29"
$pwm.Range("F29").Formula = "=C29+2*D29"
$pwm.Range("G29").Formula = "=E29/F29"
$pwm.Range("G29").Style = "Normal"

# Stop row
$pwm.Range("B30").Value = "Stop"
$pwm.Range("C30").Value = 4700
$pwm.Range("E30").Formula = "=C30+D30"
$pwm.Range("F30").Formula = "=C30+2*D30"
$pwm.Range("G30").Formula = "=E30/F30"
$pwm.Range("G30").Style = "Normal"

# Apply the "#,##0" number format to the cells used for raw R1 and Ton totals
$pwm.Range("C27").NumberFormat = "#,##0"
$pwm.Range("E27").NumberFormat = "#,##0"
$pwm.Range("E28").NumberFormat = "#,##0"
$pwm.Range("E29").NumberFormat = "#,##0"
$pwm.Range("E30").NumberFormat = "#,##0"

# Blank "touched" cells surrounding the new table (mirrors the original edit)
$pwm.Range("G25").Style = "Normal"
$pwm.Range("G31").Style = "Normal"

$pwm.Range("G14").Select()
